$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$forceTextCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D15", "D17", "D20", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D41", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.890.99'
$ws.Range("E2").Value = '  -5.27%  '
$ws.Range("D3").Value = '2.205.86'
$ws.Range("E3").Value = '  -6.79%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '313.23'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = '97.38'
$ws.Range("E6").Value = '  -9.91%  '
$ws.Range("E7").Value = '  -7.67%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.554'
$ws.Range("E9").Value = '  -10.03%  '
$ws.Range("D10").Value = '36.33'
$ws.Range("E10").Value = '  -11.58%  '
$ws.Range("D11").Value = '54.13'
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("D12").Value = '0.0822'
$ws.Range("E12").Value = '  -10.57%  '
$ws.Range("D13").Value = '7.69'
$ws.Range("E13").Value = '  -9.28%  '
$ws.Range("E14").Value = '  -4.06%  '
$ws.Range("D15").Value = '0.859'
$ws.Range("E15").Value = '  -12.32%  '
$ws.Range("D16").Value = '2.541.99'
$ws.Range("E16").Value = '  -6.77%  '
$ws.Range("D17").Value = '14.02'
$ws.Range("E17").Value = '  -7.98%  '
$ws.Range("D18").Value = '2.204.33'
$ws.Range("E18").Value = '  -7.05%  '
$ws.Range("D19").Value = '42.734.51'
$ws.Range("E19").Value = '  -5.55%  '
$ws.Range("D20").Value = '14.63'
$ws.Range("E20").Value = '  -1.53%  '
$ws.Range("D21").Value = '0.0₃0952'
$ws.Range("E21").Value = '  -10.35%  '
$ws.Range("D22").Value = '6.35'
$ws.Range("E22").Value = '  -13.05%  '
$ws.Range("D23").Value = '65.09'
$ws.Range("E23").Value = '  -11.06%  '
$ws.Range("D24").Value = '3.14'
$ws.Range("E24").Value = '  -9.72%  '
$ws.Range("D25").Value = '235.42'
$ws.Range("E25").Value = '  -9.62%  '
$ws.Range("E26").Value = '  -8.41%  '
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  -10.29%  '
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  -2.61%  '
$ws.Range("D30").Value = '6.18'
$ws.Range("E30").Value = '  -15.33%  '
$ws.Range("D31").Value = '20.37'
$ws.Range("E31").Value = '  -8.83%  '
$ws.Range("D32").Value = '0.0873'
$ws.Range("E32").Value = '  -9.50%  '
$ws.Range("D33").Value = '33.57'
$ws.Range("E33").Value = '  -10.96%  '
$ws.Range("D34").Value = '154.26'
$ws.Range("E34").Value = '  -8.69%  '
$ws.Range("E35").Value = '  -6.20%  '
$ws.Range("D36").Value = '3.14'
$ws.Range("E36").Value = '  +5.76%  '
$ws.Range("D37").Value = '1.97'
$ws.Range("E37").Value = '  +12.85%  '
$ws.Range("E38").Value = '  -6.70%  '
$ws.Range("D39").Value = '4.39'
$ws.Range("E39").Value = '  -8.13%  '
$ws.Range("E40").Value = '  -12.79%  '
$ws.Range("D41").Value = '3.67'
$ws.Range("E41").Value = '  -6.55%  '
$ws.Range("E42").Value = '  -8.91%  '
$ws.Range("D43").Value = '1.859.49'
$ws.Range("E43").Value = '  +9.22%  '
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '12.10'
$ws.Range("E45").Value = '  -6.28%  '
$ws.Range("D46").Value = '88.47'
$ws.Range("E46").Value = '  -11.56%  '
$ws.Range("D47").Value = '0.205'
$ws.Range("E47").Value = '  -11.10%  '
$ws.Range("D48").Value = '5.39'
$ws.Range("E48").Value = '  -2.67%  '
$ws.Range("D49").Value = '75.51'
$ws.Range("E49").Value = '  -6.54%  '
$ws.Range("E50").Value = '  -14.04%  '
$ws.Range("D51").Value = '8.60'
$ws.Range("E51").Value = '  -6.11%  '

foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
